# Update time-series factor values for fertility and cohabitation alignment parameters
# (cohabitation adjustment series in UK_cohabitation_adjustment, fertility adjustment series
# in UK_fertility_adjustment). Also extends both series through 2070 (row 62).

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("UK_cohabitation_adjustment")
$ws4 = $wb.Worksheets.Item("UK_fertility_adjustment")

# UK_cohabitation_adjustment (fka sheet3)
$ws3.Cells.Item(4, 2).Value = -0.95768509819680903
$ws3.Cells.Item(5, 2).Value = -0.83295781453007001
$ws3.Cells.Item(6, 2).Value = -0.80485660959326355
$ws3.Cells.Item(7, 2).Value = -0.78981867057924415
$ws3.Cells.Item(8, 2).Value = -0.77794475503792282
$ws3.Cells.Item(9, 2).Value = -0.75171647264083308
$ws3.Cells.Item(10, 2).Value = -0.74109647439251236
$ws3.Cells.Item(11, 2).Value = -0.72776510107134573
$ws3.Cells.Item(12, 2).Value = -0.7147060419149136
$ws3.Cells.Item(13, 2).Value = -0.70240616614065376
$ws3.Cells.Item(14, 2).Value = -0.69105367590666433
$ws3.Cells.Item(15, 2).Value = -0.67898522305967135
$ws3.Cells.Item(16, 2).Value = -0.67837200227633576
$ws3.Cells.Item(17, 2).Value = -0.66865521484219381
$ws3.Cells.Item(18, 2).Value = -0.66213777125088935
$ws3.Cells.Item(19, 2).Value = -0.65624918555298262
$ws3.Cells.Item(20, 2).Value = -0.65247439548232655
$ws3.Cells.Item(21, 2).Value = -0.64776174741668557
$ws3.Cells.Item(22, 2).Value = -0.64184148070049241
$ws3.Cells.Item(23, 2).Value = -0.63931566276826524
$ws3.Cells.Item(24, 2).Value = -0.63638760354965007
$ws3.Cells.Item(25, 2).Value = -0.63218797948615213
$ws3.Cells.Item(26, 2).Value = -0.62979949109872213
$ws3.Cells.Item(27, 2).Value = -0.62642840586058846
$ws3.Cells.Item(28, 2).Value = -0.62703195252281929
$ws3.Cells.Item(29, 2).Value = -0.62763592693798487
$ws3.Cells.Item(30, 2).Value = -0.62901387736233438
$ws3.Cells.Item(31, 2).Value = -0.62870490298773551
$ws3.Cells.Item(32, 2).Value = -0.63132275325212173
$ws3.Cells.Item(33, 2).Value = -0.63189417814259741
$ws3.Cells.Item(34, 2).Value = -0.63431628626780423
$ws3.Cells.Item(35, 2).Value = -0.63319756552905415
$ws3.Cells.Item(36, 2).Value = -0.63471667603087323
$ws3.Cells.Item(37, 2).Value = -0.63455648067220116
$ws3.Cells.Item(38, 2).Value = -0.63527012935003857
$ws3.Cells.Item(39, 2).Value = -0.63488663198137396
$ws3.Cells.Item(40, 2).Value = -0.6327623871436846
$ws3.Cells.Item(41, 2).Value = -0.63172645653828041
$ws3.Cells.Item(42, 2).Value = -0.63131309782462763
$ws3.Cells.Item(43, 2).Value = -0.62711212192261312
$ws3.Cells.Item(44, 2).Value = -0.62483324297863929
$ws3.Cells.Item(45, 2).Value = -0.62289315139452228
$ws3.Cells.Item(46, 2).Value = -0.6208575795852076
$ws3.Cells.Item(47, 2).Value = -0.62046331508145192
$ws3.Cells.Item(48, 2).Value = -0.61860426418741976
$ws3.Cells.Item(49, 2).Value = -0.61533367788355553
$ws3.Cells.Item(50, 2).Value = -0.61489963029454253
$ws3.Cells.Item(51, 2).Value = -0.61364330813003531
$ws3.Cells.Item(52, 2).Value = -0.61242239408296961
$ws3.Cells.Item(53, 2).Value = -0.61091556383364465
$ws3.Cells.Item(54, 2).Value = -0.60932285274850906
$ws3.Cells.Item(55, 2).Value = -0.60791630071486702
$ws3.Cells.Item(56, 2).Value = -0.60902347078757146
$ws3.Cells.Item(57, 2).Value = -0.60734177454365323
$ws3.Cells.Item(58, 2).Value = -0.60672079418846891
$ws3.Cells.Item(59, 2).Value = -0.6058496118125366
$ws3.Cells.Item(60, 2).Value = -0.60652410750329455
$ws3.Cells.Item(61, 2).Value = -0.60652410750329455
$ws3.Cells.Item(62, 1).Value = 2070
$ws3.Cells.Item(62, 2).Value = -0.60652410750329455

# UK_fertility_adjustment (fka sheet4)
$ws4.Cells.Item(4, 2).Value = -0.48463199296432402
$ws4.Cells.Item(5, 2).Value = -0.44279476957698699
$ws4.Cells.Item(6, 2).Value = -0.39501557684159466
$ws4.Cells.Item(7, 2).Value = -0.37898938142570604
$ws4.Cells.Item(8, 2).Value = -0.35746845133604671
$ws4.Cells.Item(9, 2).Value = -0.33551640350223444
$ws4.Cells.Item(10, 2).Value = -0.31611792786789245
$ws4.Cells.Item(11, 2).Value = -0.29419113555012588
$ws4.Cells.Item(12, 2).Value = -0.275174972273524
$ws4.Cells.Item(13, 2).Value = -0.25957357188884783
$ws4.Cells.Item(14, 2).Value = -0.2522555041315423
$ws4.Cells.Item(15, 2).Value = -0.2484037751560034
$ws4.Cells.Item(16, 2).Value = -0.24178593173429327
$ws4.Cells.Item(17, 2).Value = -0.24449858212979714
$ws4.Cells.Item(18, 2).Value = -0.24903198566118756
$ws4.Cells.Item(19, 2).Value = -0.25264201070565656
$ws4.Cells.Item(20, 2).Value = -0.25364233752894055
$ws4.Cells.Item(21, 2).Value = -0.24768152163486329
$ws4.Cells.Item(22, 2).Value = -0.24974442678892514
$ws4.Cells.Item(23, 2).Value = -0.24926794764317603
$ws4.Cells.Item(24, 2).Value = -0.250134489952685
$ws4.Cells.Item(25, 2).Value = -0.25159703697288732
$ws4.Cells.Item(26, 2).Value = -0.25383394104860174
$ws4.Cells.Item(27, 2).Value = -0.26016905205776741
$ws4.Cells.Item(28, 2).Value = -0.27274303604096028
$ws4.Cells.Item(29, 2).Value = -0.27650991999973468
$ws4.Cells.Item(30, 2).Value = -0.28460685361503957
$ws4.Cells.Item(31, 2).Value = -0.29665890398267214
$ws4.Cells.Item(32, 2).Value = -0.30980188741828268
$ws4.Cells.Item(33, 2).Value = -0.3249147353451094
$ws4.Cells.Item(34, 2).Value = -0.33013144629608343
$ws4.Cells.Item(35, 2).Value = -0.33400456174430643
$ws4.Cells.Item(36, 2).Value = -0.34264176458383305
$ws4.Cells.Item(37, 2).Value = -0.34821267929122846
$ws4.Cells.Item(38, 2).Value = -0.34456161150133158
$ws4.Cells.Item(39, 2).Value = -0.34561532498119613
$ws4.Cells.Item(40, 2).Value = -0.33873736519970848
$ws4.Cells.Item(41, 2).Value = -0.33739651765499096
$ws4.Cells.Item(42, 2).Value = -0.33804860844064244
$ws4.Cells.Item(43, 2).Value = -0.33657044977191414
$ws4.Cells.Item(44, 2).Value = -0.33656675024980498
$ws4.Cells.Item(45, 2).Value = -0.33985483078933226
$ws4.Cells.Item(46, 2).Value = -0.33668670713123178
$ws4.Cells.Item(47, 2).Value = -0.33750504731956005
$ws4.Cells.Item(48, 2).Value = -0.34079779337557842
$ws4.Cells.Item(49, 2).Value = -0.33872371122332023
$ws4.Cells.Item(50, 2).Value = -0.33869137298543883
$ws4.Cells.Item(51, 2).Value = -0.34103219793504497
$ws4.Cells.Item(52, 2).Value = -0.34050601889900228
$ws4.Cells.Item(53, 2).Value = -0.33767373591580502
$ws4.Cells.Item(54, 2).Value = -0.33668796343699775
$ws4.Cells.Item(55, 2).Value = -0.33369114324283816
$ws4.Cells.Item(56, 2).Value = -0.33282968867461227
$ws4.Cells.Item(57, 2).Value = -0.32675564471673357
$ws4.Cells.Item(58, 2).Value = -0.32143826156637562
$ws4.Cells.Item(59, 2).Value = -0.31711471655916201
$ws4.Cells.Item(60, 2).Value = -0.31670999077073719
$ws4.Cells.Item(61, 2).Value = -0.31670999077073719
$ws4.Cells.Item(62, 1).Value = 2070
$ws4.Cells.Item(62, 2).Value = -0.31670999077073719

# Tidy up the sheet views: scroll back to the top-left and drop any
# leftover multi-cell selection on both updated sheets. Re-activate the
# fertility-adjustment sheet last so it remains the active tab, matching
# the workbook's original activeTab.
$ws3.Activate() | Out-Null
$ws3.Range("A1").Select() | Out-Null
$ws4.Activate() | Out-Null
$ws4.Range("A1").Select() | Out-Null
$ws4.Activate() | Out-Null
